# Update "想去人数" (interest / want-to-go counter) values that changed
# between the previous gh-pages data snapshot and the one generated at
# commit 456a3b4.
#
# Sheet order in the workbook:
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活  (Local life)
#   4 = 全部类型  (All types)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 9440
$ws.Range("F7").Value = 6469
$ws.Range("F10").Value = 9898
$ws.Range("F11").Value = 11346
$ws.Range("F13").Value = 1168
$ws.Range("F14").Value = 4967
$ws.Range("F15").Value = 806
$ws.Range("F21").Value = 263
$ws.Range("F23").Value = 897
$ws.Range("F24").Value = 1272
$ws.Range("F29").Value = 636
$ws.Range("F31").Value = 190
$ws.Range("F32").Value = 1781
$ws.Range("F33").Value = 95
$ws.Range("F38").Value = 38
$ws.Range("F41").Value = 88
$ws.Range("F48").Value = 4219

# --- Sheet 2: 演出 -----------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Range("F19").Value = 11

# --- Sheet 3: 本地生活 -------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5976

# --- Sheet 4: 全部类型 -------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 9440
$ws.Range("F8").Value = 6469
$ws.Range("F10").Value = 9898
$ws.Range("F11").Value = 11346
$ws.Range("F13").Value = 1168
$ws.Range("F14").Value = 4967
$ws.Range("F15").Value = 806
$ws.Range("F22").Value = 263
$ws.Range("F28").Value = 636
$ws.Range("F30").Value = 190
$ws.Range("F31").Value = 1781
$ws.Range("F32").Value = 95
$ws.Range("F42").Value = 38
$ws.Range("F49").Value = 4219

Write-Host "Applied 33 cell updates across 4 sheets"
